$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; existing rows 20..82 shift down to 21..83,
# and formatting (e.g. the date style on column D) is inherited from the
# row that was previously there.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record.
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(20, 4).Value = 44715
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = 100112021
$ws.Cells.Item(20, 7).Value = "Ají"
$ws.Cells.Item(20, 8).Value = "Inferno"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 120
$ws.Cells.Item(20, 11).Value = 17000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 13).Value = 17500
$ws.Cells.Item(20, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(20, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(20, 16).Value = 1167
$ws.Cells.Item(20, 17).Value = 15
$ws.Cells.Item(20, 18).Value = "Hortaliza"
